$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 216; this shifts the existing rows 216..259
# down to 217..260 (carrying their data/formatting with them), matching
# the diff's observed "every row from 216 on is shifted down by one" plus
# a brand new row 260 appearing (built from the tail of the old data).
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new data record.
$ws.Cells.Item(216, 1).Value2 = 6
$ws.Cells.Item(216, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(216, 3).Value2 = "Metropolitana"
$ws.Cells.Item(216, 4).Value2 = 44889
$ws.Cells.Item(216, 5).Value2 = 13
$ws.Cells.Item(216, 6).Value2 = 100112029
$ws.Cells.Item(216, 7).Value2 = "Orégano"
$ws.Cells.Item(216, 8).Value2 = "Sin especificar"
$ws.Cells.Item(216, 9).Value2 = "Primera"
$ws.Cells.Item(216, 10).Value2 = 47
$ws.Cells.Item(216, 11).Value2 = 16000
$ws.Cells.Item(216, 12).Value2 = 17000
$ws.Cells.Item(216, 13).Value2 = 16447
$ws.Cells.Item(216, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(216, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(216, 16).Value2 = 5482
$ws.Cells.Item(216, 17).Value2 = 3
$ws.Cells.Item(216, 18).Value2 = "Hortaliza"
